$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 47/48: Cronos <-> HuobiToken swap positions (with updated price/volume)
$ws.Range("B47").Value = "HuobiToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"

# Price (column D) updates
$ws.Range("D2").Value = "36.753.91"
$ws.Range("D3").Value = "1.924.55"
$ws.Range("D4").Formula = "'1.00"
$ws.Range("D5").Formula = "'251.23"
$ws.Range("D6").Formula = "'0.699"
$ws.Range("D7").Formula = "'1.00"
$ws.Range("D8").Formula = "'44.28"
$ws.Range("D9").Formula = "'59.11"
$ws.Range("D12").Formula = "'0.0995"
$ws.Range("D13").Formula = "'14.49"
$ws.Range("D14").Formula = "'0.802"
$ws.Range("D15").Value = "2.204.11"
$ws.Range("D17").Value = "1.929.54"
$ws.Range("D18").Value = "36.648.39"
$ws.Range("D19").Formula = "'74.52"
$ws.Range("D20").Value = "0.0₃0867"
$ws.Range("D21").Formula = "'251.68"
$ws.Range("D22").Formula = "'13.31"
$ws.Range("D23").Formula = "'5.23"
$ws.Range("D24").Formula = "'2.68"
$ws.Range("D26").Formula = "'2.23"
$ws.Range("D27").Formula = "'168.08"
$ws.Range("D28").Formula = "'8.80"
$ws.Range("D29").Formula = "'18.83"
$ws.Range("D30").Formula = "'0.129"
$ws.Range("D32").Formula = "'0.0610"
$ws.Range("D33").Formula = "'2.00"
$ws.Range("D36").Formula = "'0.0859"
$ws.Range("D39").Formula = "'17.73"
$ws.Range("D40").Formula = "'2.02"
$ws.Range("D41").Formula = "'109.90"
$ws.Range("D43").Formula = "'17.22"
$ws.Range("D45").Value = "1.346.40"
$ws.Range("D46").Formula = "'2.37"
$ws.Range("D47").Formula = "'2.49"
$ws.Range("D48").Formula = "'0.0813"
$ws.Range("D50").Formula = "'6.45"
$ws.Range("D51").Value = "2.111.05"

# Volume(1h) percentage (column E) updates
$ws.Range("E2").Value = "  +3.91%  "
$ws.Range("E3").Value = "  +2.10%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("E5").Value = "  +2.31%  "
$ws.Range("E6").Value = "  +1.53%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +1.40%  "
$ws.Range("E9").Value = "  +10.57%  "
$ws.Range("E10").Value = "  +3.71%  "
$ws.Range("E12").Value = "  +2.46%  "
$ws.Range("E14").Value = "  +4.99%  "
$ws.Range("E15").Value = "  +2.12%  "
$ws.Range("E16").Value = "  +4.82%  "
$ws.Range("E17").Value = "  +2.54%  "
$ws.Range("E18").Value = "  +3.15%  "
$ws.Range("E19").Value = "  +2.00%  "
$ws.Range("E20").Value = "  +5.43%  "
$ws.Range("E21").Value = "  +2.89%  "
$ws.Range("E22").Value = "  +3.89%  "
$ws.Range("E23").Value = "  +5.54%  "
$ws.Range("E24").Value = "  +1.20%  "
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("E26").Value = "  +3.49%  "
$ws.Range("E27").Value = "  +1.42%  "
$ws.Range("E28").Value = "  +3.09%  "
$ws.Range("E29").Value = "  +2.72%  "
$ws.Range("E30").Value = "  +1.78%  "
$ws.Range("E31").Value = "  +6.30%  "
$ws.Range("E32").Value = "  +3.76%  "
$ws.Range("E33").Value = "  +4.55%  "
$ws.Range("E34").Value = "  +4.72%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("E36").Value = "  +22.73%  "
$ws.Range("E37").Value = "  -13.84%  "
$ws.Range("E38").Value = "  +2.34%  "
$ws.Range("E39").Value = "  +45.91%  "
$ws.Range("E40").Value = "  +3.53%  "
$ws.Range("E41").Value = "  +14.18%  "
$ws.Range("E42").Value = "  +5.08%  "
$ws.Range("E43").Value = "  -0.80%  "
$ws.Range("E44").Value = "  +3.45%  "
$ws.Range("E45").Value = "  +3.03%  "
$ws.Range("E46").Value = "  +1.44%  "
$ws.Range("E47").Value = "  +4.81%  "
$ws.Range("E48").Value = "  +1.93%  "
$ws.Range("E49").Value = "  +3.11%  "
$ws.Range("E50").Value = "  +3.22%  "
$ws.Range("E51").Value = "  +2.26%  "
